$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "teste"
$ws.Range("B1").Value = "A/C"
$ws.Range("C1").Value = 23
$ws.Range("D1").Value = $false

$ws.Range("A2").Value = "teste"
$ws.Range("B2").Value = "A/C"
$ws.Range("C2").Value = 23
$ws.Range("D2").Value = $false

$ws.Range("A3").Value = "new ar"
$ws.Range("B3").Value = "A/C"
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = $false

$ws.Range("A4").Value = "lamp1"
$ws.Range("B4").Value = "Lâmpada"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = $false
